$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old emoji values in column A (statut) to their new replacement text
$map = @{
    "📘" = "⚠️"
    "📕" = "-3"
    "📙" = "+3"
    "📗" = "✅"
}

# These replacement values look like numbers, so Excel would silently convert
# them to numeric cells unless we force the cell to stay text.
$numericLooking = @{ "-3" = $true; "+3" = $true }

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $newVal = $map[$val]
        if ($numericLooking.ContainsKey($newVal)) {
            # Temporarily force text formatting so "-3"/"+3" are stored as
            # shared-string text instead of being coerced into numbers, then
            # restore the default "Normal" style so the cell keeps looking
            # like the rest of the column.
            $cell.NumberFormat = "@"
            $cell.Value = $newVal
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newVal
        }
    }
}
